# rfctr: isolate Slides.add_slide() to single test
# Add a second slide (a table slide) to the deck and tweak slide 1's
# textbox (resize + add a hyperlinked "yahoo.com" paragraph).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 1: resize the textbox and append a hyperlinked paragraph
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$shp = $s1.Shapes.Item(1)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Resize the textbox to its new target size (EMU -> points, 12700 EMU/pt)
$shp.Width = 1226455 / 12700
$shp.Height = 1200329 / 12700

# Append a new paragraph "yahoo.com" after the existing text, then turn
# just that new run into a hyperlink.
$cr = [char]13
$tr.InsertAfter($cr + "yahoo.com")
$fullText = $tf.TextRange.Text
$linkText = "yahoo.com"
$startPos = $fullText.Length - $linkText.Length + 1
$linkRange = $tf.TextRange.Characters($startPos, $linkText.Length)
$linkRange.ActionSettings.Item(1).Hyperlink.Address = "http://www.yahoo.com/"

# ---------------------------------------------------------------
# Slide 2 (new): a blank slide holding a 2x2 table
# ---------------------------------------------------------------
$s2 = $p.Slides.Add(2, 12)
$tbl = $s2.Shapes.AddTable(2, 2, 1524000 / 12700, 1397000 / 12700, 6096000 / 12700, 741680 / 12700)
$tbl.Table.Cell(1, 1).Shape.TextFrame.TextRange.Text = "Text run in table cell"
